# Rename the "patch version" suffixes used in column A of the QueryResult
# sheet: "-Repaired" rows become "-Auto", and "-Fixed" rows become "-Manual".
# ("-Buggy" rows are left untouched.)
#
# The two header rows (A18 / A220, literal text "201 projects") and the
# summary rows (Sum/Average/Minimum/Maximum/Standard deviation/Variance)
# are untouched — their shared-string slot simply shifts because of the
# dedup/removal of the old "-Fixed"/"-Repaired" strings, which Excel
# handles automatically when the cell values are rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Process the "-Repaired" -> "-Auto" block BEFORE the "-Fixed" -> "-Manual"
# block so that newly introduced shared strings keep the same relative
# ordering as produced by the original authoring tool (Auto before Manual).
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -ne $null -and $val -is [string] -and $val.EndsWith("-Repaired")) {
        $base = $val.Substring(0, $val.Length - "-Repaired".Length)
        $cell.Value = $base + "-Auto"
    }
}

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -ne $null -and $val -is [string] -and $val.EndsWith("-Fixed")) {
        $base = $val.Substring(0, $val.Length - "-Fixed".Length)
        $cell.Value = $base + "-Manual"
    }
}
